# Append the new anime-character record as row 3 (mirrors row 2's layout):
# id, nombre, nombre_kanji, role, anime_id, imagen_url, fecha_extraccion
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = ""
$ws.Cells.Item(3, 3).Value = ""
$ws.Cells.Item(3, 4).Value = "Main"
$ws.Cells.Item(3, 5).Value = 853
$ws.Cells.Item(3, 6).Value = ""
$ws.Cells.Item(3, 7).Value = "2025-03-02 15:25:53"
